$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.963.26'
$ws.Range('E2').Value = '  -3.92%  '

$ws.Range('D3').Value = '1.641.00'
$ws.Range('E3').Value = '  -5.78%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9987'
$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.00'
$ws.Range('E5').Value = '  -5.77%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.0000'
$ws.Range('E6').Value = '  -0.04%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4756'
$ws.Range('E7').Value = '  -5.26%  '

$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '39.55'
$ws.Range('E8').Value = '  -3.47%  '

$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2585'
$ws.Range('E9').Value = '  -5.66%  '

$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06083'
$ws.Range('E10').Value = '  -1.81%  '

$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07038'
$ws.Range('E11').Value = '  -3.07%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.650.71'
$ws.Range('E12').Value = '  -5.25%  '

$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.56'
$ws.Range('E13').Value = '  -4.12%  '

$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5892'
$ws.Range('E14').Value = '  -9.98%  '

$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.325'
$ws.Range('E15').Value = '  -8.16%  '

$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '73.68'
$ws.Range('E16').Value = '  -5.22%  '

$ws.Range('B17').Value = 'Dai'
$ws.Range('C17').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.000'
$ws.Range('E17').Value = '  -0.03%  '

$ws.Range('B18').Value = 'BinanceUSD'
$ws.Range('C18').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9997'
$ws.Range('E18').Value = '  -0.06%  '

$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '24.950.80'
$ws.Range('E19').Value = '  -4.04%  '

$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000006571'
$ws.Range('E20').Value = '  -4.16%  '

$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.22'
$ws.Range('E21').Value = '  -5.70%  '

$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '1.860.39'
$ws.Range('E22').Value = '  -5.48%  '

$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.310'
$ws.Range('E23').Value = '  -6.66%  '

$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.540'
$ws.Range('E24').Value = '  -2.48%  '

$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.236'
$ws.Range('E25').Value = '  -3.23%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '133.75'
$ws.Range('E26').Value = '  -0.89%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.90'
$ws.Range('E27').Value = '  -2.45%  '

$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.379'
$ws.Range('E28').Value = '  -8.01%  '

$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '104.26'
$ws.Range('E29').Value = '  -1.18%  '

$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.639'
$ws.Range('E30').Value = '  -8.50%  '

$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.890'
$ws.Range('E31').Value = '  -2.19%  '

$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.07595'
$ws.Range('E32').Value = '  -6.97%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.571'
$ws.Range('E33').Value = '  -3.54%  '

$ws.Range('B34').Value = 'Frax'
$ws.Range('C34').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9995'
$ws.Range('E34').Value = '  +0.00%  '

$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.04270'
$ws.Range('E35').Value = '  -10.08%  '

$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.573'
$ws.Range('E36').Value = '  -3.60%  '

$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5933'
$ws.Range('E37').Value = '  -3.29%  '

$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9278'
$ws.Range('E38').Value = '  -7.18%  '

$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.570'
$ws.Range('E39').Value = '  -6.78%  '

$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8619'
$ws.Range('E40').Value = '  +7.42%  '

$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9998'
$ws.Range('E41').Value = '  -0.03%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01495'
$ws.Range('E42').Value = '  -7.94%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '98.65'
$ws.Range('E43').Value = '  -2.56%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.766'
$ws.Range('E44').Value = '  -8.91%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3700'
$ws.Range('E45').Value = '  -5.71%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.667'
$ws.Range('E46').Value = '  -7.58%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1101'
$ws.Range('E47').Value = '  -6.61%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.112'
$ws.Range('E48').Value = '  -4.54%  '

$ws.Range('E49').Value = '  -1.80%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.001'
$ws.Range('E50').Value = '  -0.08%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '28.74'
$ws.Range('E51').Value = '  -7.19%  '
